$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.030.45'

$ws.Range("D3").Value = '3.311.70'
$ws.Range("E3").Value = '  -0.86%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = "'" + '585.34'
$ws.Range("E5").Value = '  +1.98%  '

$ws.Range("D6").Value = "'" + '182.33'
$ws.Range("E6").Value = '  +0.90%  '

$ws.Range("E7").Value = '  +2.55%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").Value = '3.306.33'
$ws.Range("E9").Value = '  -0.96%  '

$ws.Range("E10").Value = '  -3.52%  '

$ws.Range("D11").Value = "'" + '6.81'
$ws.Range("E11").Value = '  +2.55%  '

$ws.Range("D12").Value = "'" + '0.400'
$ws.Range("E12").Value = '  -0.88%  '

$ws.Range("D13").Value = '3.885.51'
$ws.Range("E13").Value = '  -0.79%  '

$ws.Range("E14").Value = '  -3.06%  '

$ws.Range("D15").Value = '66.075.17'
$ws.Range("E15").Value = '  -0.97%  '

$ws.Range("D16").Value = "'" + '26.12'
$ws.Range("E16").Value = '  -3.43%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.346.48'
$ws.Range("E17").Value = '  +0.34%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = "'" + '0.0000165'
$ws.Range("E18").Value = '  -1.72%  '

$ws.Range("D19").Value = "'" + '424.03'
$ws.Range("E19").Value = '  -2.71%  '

$ws.Range("E20").Value = '  -2.69%  '

$ws.Range("D21").Value = "'" + '13.12'
$ws.Range("E21").Value = '  -3.54%  '

$ws.Range("E22").Value = '  -2.99%  '

$ws.Range("E23").Value = '  -2.60%  '

$ws.Range("E24").Value = '  +0.09%  '

$ws.Range("E25").Value = '  +0.17%  '

$ws.Range("D26").Value = '3.462.75'
$ws.Range("E26").Value = '  -0.63%  '

$ws.Range("E28").Value = '  +5.12%  '

$ws.Range("E29").Value = '  -4.09%  '

$ws.Range("D30").Value = "'" + '8.88'
$ws.Range("E30").Value = '  -1.49%  '

$ws.Range("D31").Value = "'" + '1.00'
$ws.Range("E31").Value = '  +0.46%  '

$ws.Range("D32").Value = "'" + '1.90'
$ws.Range("E32").Value = '  -3.01%  '

$ws.Range("D33").Value = "'" + '22.34'
$ws.Range("E33").Value = '  -2.28%  '

$ws.Range("D34").Value = "'" + '1.00'

$ws.Range("D35").Value = "'" + '5.15'

$ws.Range("D36").Value = "'" + '6.53'
$ws.Range("E36").Value = '  -3.37%  '

$ws.Range("E37").Value = '  -5.01%  '

$ws.Range("D38").Value = "'" + '160.45'
$ws.Range("E38").Value = '  -1.49%  '

$ws.Range("E39").Value = '  -3.77%  '

$ws.Range("D40").Value = '2.871.22'
$ws.Range("E40").Value = '  +1.63%  '

$ws.Range("E41").Value = '  -1.65%  '

$ws.Range("D42").Value = "'" + '26.23'
$ws.Range("E42").Value = '  -4.56%  '

$ws.Range("D43").Value = "'" + '0.758'
$ws.Range("E43").Value = '  -4.49%  '

$ws.Range("D44").Value = "'" + '4.29'
$ws.Range("E44").Value = '  -2.89%  '

$ws.Range("D45").Value = "'" + '39.88'
$ws.Range("E45").Value = '  -0.78%  '

$ws.Range("D46").Value = "'" + '0.0658'
$ws.Range("E46").Value = '  -1.17%  '

$ws.Range("D47").Value = "'" + '5.90'
$ws.Range("E47").Value = '  -4.79%  '

$ws.Range("E48").Value = '  -3.42%  '

$ws.Range("D49").Value = "'" + '313.19'
$ws.Range("E49").Value = '  -2.51%  '

$ws.Range("D50").Value = "'" + '23.04'
$ws.Range("E50").Value = '  -5.78%  '

$ws.Range("D51").Value = "'" + '0.0271'
$ws.Range("E51").Value = '  -1.16%  '
